# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted in
# front of the "Late" column, pushing "Late" / "heading" ("Date") /
# "Outstanding" one column to the right (N -> O, O -> P, P -> Q).
# The new column takes on the same look/width as column M ("In Advance").
#
# Also, the "Repayment schedule" sheet becomes the active/selected sheet
# of the workbook (it previously was "Acc_Repayment").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# remaining columns (Late / heading / Outstanding) one place to the right.
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as column M ("In Advance").
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet / tab, with cell R6 selected
# (this also clears the previous selection/active state on whichever sheet
# used to be active).
$ws.Activate()
$ws.Range("R6").Select()
